# Insert a new data row at row 54 (pushing rows 54..134 down to 55..135),
# then populate the newly inserted row with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a whole row above the current row 54 - this shifts rows 54-134
# down to 55-135 and carries the row-53 formatting (incl. date style) onto
# the new row, matching how column D's other cells are formatted.
$ws.Rows.Item(54).Insert()

# Fill the newly-inserted row 54 with the new record's values.
$ws.Cells.Item(54, 1).Value = 11
$ws.Cells.Item(54, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(54, 3).Value = "Bíobío"
$ws.Cells.Item(54, 4).Value = 44579
$ws.Cells.Item(54, 5).Value = 8
$ws.Cells.Item(54, 6).Value = "Fruta"
$ws.Cells.Item(54, 7).Value = 100108
$ws.Cells.Item(54, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(54, 9).Value = 100108005
$ws.Cells.Item(54, 10).Value = "Piña"
$ws.Cells.Item(54, 11).Value = "Caramelo"
$ws.Cells.Item(54, 12).Value = "Segunda"
$ws.Cells.Item(54, 13).Value = 250
$ws.Cells.Item(54, 14).Value = 14000
$ws.Cells.Item(54, 15).Value = 15000
$ws.Cells.Item(54, 16).Value = 14520
$ws.Cells.Item(54, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(54, 18).Value = "Ecuador"
$ws.Cells.Item(54, 19).Value = 1037
$ws.Cells.Item(54, 20).Value = 14
